$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.073327139394411
$ws.Range("D2").Value = 1.075973923348161
$ws.Range("E2").Value = 1.076538476814783
$ws.Range("F2").Value = 1.085503883737039
$ws.Range("I2").Value = 1.05183948584513
$ws.Range("J2").Value = 1.078242206291264
$ws.Range("K2").Value = 1.078659079184975
$ws.Range("L2").Value = 1.079222144536872
$ws.Range("M2").Value = 1.08816414337042
$ws.Range("N2").Value = 1.079773433900307
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.075138003723011
$ws.Range("D3").Value = 1.077703532402589
$ws.Range("E3").Value = 1.078168336174205
$ws.Range("F3").Value = 1.087214261235087
$ws.Range("I3").Value = 1.052351135115382
$ws.Range("J3").Value = 1.079707697141129
$ws.Range("K3").Value = 1.080203873611861
$ws.Range("L3").Value = 1.080667542307511
$ws.Range("M3").Value = 1.08969158436577
$ws.Range("N3").Value = 1.081241005915274
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.076306281046743
$ws.Range("D4").Value = 1.078819539352444
$ws.Range("E4").Value = 1.079219552491735
$ws.Range("F4").Value = 1.088317828673072
$ws.Range("I4").Value = 1.05267907201375
$ws.Range("J4").Value = 1.080652181150632
$ws.Range("K4").Value = 1.081199831930304
$ws.Range("L4").Value = 1.08159891818876
$ws.Range("M4").Value = 1.090676309610072
$ws.Range("N4").Value = 1.082186831200404
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.076796611188501
$ws.Range("D5").Value = 1.079287967057067
$ws.Range("E5").Value = 1.079660682838797
$ws.Range("F5").Value = 1.088781027640585
$ws.Range("I5").Value = 1.052816192296552
$ws.Range("J5").Value = 1.081048350408905
$ws.Range("K5").Value = 1.081617679694437
$ws.Range("L5").Value = 1.081989551377564
$ws.Range("M5").Value = 1.091089433773776
$ws.Range("N5").Value = 1.082583563064465
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.076878892578104
$ws.Range("D6").Value = 1.079366575110901
$ws.Range("E6").Value = 1.079734704049858
$ws.Range("F6").Value = 1.088858757807704
$ws.Range("I6").Value = 1.05283917191682
$ws.Range("J6").Value = 1.081114817097994
$ws.Range("K6").Value = 1.081687788584258
$ws.Range("L6").Value = 1.082055087048984
$ws.Range("M6").Value = 1.091158749473245
$ws.Range("N6").Value = 1.082650124143876
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.076312836034924
$ws.Range("D7").Value = 1.07882580139821
$ws.Range("E7").Value = 1.079225450022869
$ws.Range("F7").Value = 1.088324020849276
$ws.Range("I7").Value = 1.052680907138412
$ws.Range("J7").Value = 1.080657478267848
$ws.Range("K7").Value = 1.081205418564006
$ws.Range("L7").Value = 1.081604141432469
$ws.Range("M7").Value = 1.09068183313184
$ws.Range("N7").Value = 1.082192135840134
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.073939859338367
$ws.Range("D8").Value = 1.076559118159918
$ws.Range("E8").Value = 1.077090011361249
$ws.Range("F8").Value = 1.086082578042739
$ws.Range("I8").Value = 1.052013052733868
$ws.Range("J8").Value = 1.078738268704684
$ws.Range("K8").Value = 1.079181909871055
$ws.Range("L8").Value = 1.079711438522983
$ws.Range("M8").Value = 1.088681110140228
$ws.Range("N8").Value = 1.080270200779253
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.069730969466355
$ws.Range("D9").Value = 1.072539952401742
$ws.Range("E9").Value = 1.073300254874714
$ws.Range("F9").Value = 1.082107936352549
$ws.Range("I9").Value = 1.050811931359911
$ws.Range("J9").Value = 1.075326712553536
$ws.Range("K9").Value = 1.075587765101548
$ws.Range("L9").Value = 1.07634577013152
$ws.Range("M9").Value = 1.085127093293324
$ws.Range("N9").Value = 1.07685379982706
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.066905491141266
$ws.Range("D10").Value = 1.069842680790814
$ws.Range("E10").Value = 1.0707547044288
$ws.Range("F10").Value = 1.079440396178243
$ws.Range("I10").Value = 1.049994479639924
$ws.Range("J10").Value = 1.073031464981936
$ws.Range("K10").Value = 1.07317156687478
$ws.Range("L10").Value = 1.0740805573439
$ws.Range("M10").Value = 1.082737644179124
$ws.Range("N10").Value = 1.07455529274064
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.065677137699419
$ws.Range("D11").Value = 1.068670272547695
$ws.Range("E11").Value = 1.069647711505717
$ws.Range("F11").Value = 1.078280877354683
$ws.Range("I11").Value = 1.049636466104409
$ws.Range("J11").Value = 1.072032436680987
$ws.Range("K11").Value = 1.072120342395583
$ws.Range("L11").Value = 1.073094401198535
$ws.Range("M11").Value = 1.081698004602055
$ws.Range("N11").Value = 1.073554845704913
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.065220114069328
$ws.Range("D12").Value = 1.068234096026355
$ws.Range("E12").Value = 1.069235791576703
$ws.Range("F12").Value = 1.077849491199162
$ws.Range("I12").Value = 1.049502867740204
$ws.Range("J12").Value = 1.07166055807034
$ws.Range("K12").Value = 1.071729101673856
$ws.Range("L12").Value = 1.072727283793576
$ws.Range("M12").Value = 1.081311067636538
$ws.Range("N12").Value = 1.073182438983984
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.065318181866214
$ws.Range("D13").Value = 1.068327689005378
$ws.Range("E13").Value = 1.069324183305087
$ws.Range("F13").Value = 1.077942056502642
$ws.Range("I13").Value = 1.049531553048483
$ws.Range("J13").Value = 1.071740363564635
$ws.Range("K13").Value = 1.071813059245688
$ws.Range("L13").Value = 1.072806068903409
$ws.Range("M13").Value = 1.081394102006435
$ws.Range("N13").Value = 1.073262357811235
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.065639375583798
$ws.Range("D14").Value = 1.068634232291828
$ws.Range("E14").Value = 1.069613677151618
$ws.Range("F14").Value = 1.078245233015718
$ws.Range("I14").Value = 1.049625435435912
$ws.Range("J14").Value = 1.072001713402922
$ws.Range("K14").Value = 1.072088018115098
$ws.Range("L14").Value = 1.073064071893464
$ws.Range("M14").Value = 1.08166603605752
$ws.Range("N14").Value = 1.073524078796269
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.065837172469545
$ws.Range("D15").Value = 1.06882301144019
$ws.Range("E15").Value = 1.069791946088215
$ws.Range("F15").Value = 1.078431938146569
$ws.Range("I15").Value = 1.049683197595231
$ws.Range("J15").Value = 1.072162633769645
$ws.Range("K15").Value = 1.072257326865278
$ws.Range("L15").Value = 1.073222927492355
$ws.Range("M15").Value = 1.081833481205491
$ws.Range("N15").Value = 1.073685227688372
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06698690775714
$ws.Range("D16").Value = 1.06992039374219
$ws.Range("E16").Value = 1.07082807011587
$ws.Range("F16").Value = 1.079517254056228
$ws.Range("I16").Value = 1.05001815382014
$ws.Range("J16").Value = 1.073097656710871
$ws.Range("K16").Value = 1.073241226354016
$ws.Range("L16").Value = 1.074145891978257
$ws.Range("M16").Value = 1.082806534877778
$ws.Range("N16").Value = 1.074621578469421
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06770677836861
$ws.Range("D17").Value = 1.070607541351581
$ws.Range("E17").Value = 1.071476717258215
$ws.Range("F17").Value = 1.080196837121209
$ws.Range("I17").Value = 1.050227173035645
$ws.Range("J17").Value = 1.073682774753475
$ws.Range("K17").Value = 1.073857050125738
$ws.Range("L17").Value = 1.074723410370882
$ws.Range("M17").Value = 1.083415556164365
$ws.Range("N17").Value = 1.075207527446761
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068126195172117
$ws.Range("D18").Value = 1.07100791282661
$ws.Range("E18").Value = 1.071854604911215
$ws.Range("F18").Value = 1.080592797799428
$ws.Range("I18").Value = 1.050348700003999
$ws.Range("J18").Value = 1.074023566522991
$ws.Range("K18").Value = 1.074215769074791
$ws.Range("L18").Value = 1.075059756556244
$ws.Range("M18").Value = 1.083770307549264
$ws.Range("N18").Value = 1.075548803179679
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.068269126139652
$ws.Range("D19").Value = 1.071144356803671
$ws.Range("E19").Value = 1.071983377861587
$ws.Range("F19").Value = 1.080727738048838
$ws.Range("I19").Value = 1.050390071565684
$ws.Range("J19").Value = 1.074139683838872
$ws.Range("K19").Value = 1.074338002005369
$ws.Range("L19").Value = 1.075174355868215
$ws.Range("M19").Value = 1.083891187671087
$ws.Range("N19").Value = 1.075665085395467
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.067629591996764
$ws.Range("D20").Value = 1.070533861548058
$ws.Range("E20").Value = 1.071407170994198
$ws.Range("F20").Value = 1.080123968766874
$ws.Range("I20").Value = 1.050204787687302
$ws.Range("J20").Value = 1.073620048733728
$ws.Range("K20").Value = 1.073791027912952
$ws.Range("L20").Value = 1.074661501050087
$ws.Range("M20").Value = 1.083350263770283
$ws.Range("N20").Value = 1.075144712348871
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.065544813154207
$ws.Range("D21").Value = 1.068543982196474
$ws.Range("E21").Value = 1.069528448833724
$ws.Range("F21").Value = 1.07815597427505
$ws.Range("I21").Value = 1.049597806483532
$ws.Range("J21").Value = 1.071924774482744
$ws.Range("K21").Value = 1.072007070965121
$ws.Range("L21").Value = 1.072988119091838
$ws.Range("M21").Value = 1.081585979622314
$ws.Range("N21").Value = 1.073447030614
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.064229630633448
$ws.Range("D22").Value = 1.067288852839333
$ws.Range("E22").Value = 1.068342968227619
$ws.Range("F22").Value = 1.076914620364721
$ws.Range("I22").Value = 1.049212604398046
$ws.Range("J22").Value = 1.070854279392552
$ws.Range("K22").Value = 1.070880966934275
$ws.Range("L22").Value = 1.071931272302368
$ws.Range("M22").Value = 1.08047224830682
$ws.Range("N22").Value = 1.072375015297993
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.064927258158817
$ws.Range("D23").Value = 1.06795460782938
$ws.Range("E23").Value = 1.068971823647471
$ws.Range("F23").Value = 1.077573071142135
$ws.Range("I23").Value = 1.049417148133117
$ws.Range("J23").Value = 1.07142221207989
$ws.Range("K23").Value = 1.071478365017103
$ws.Range("L23").Value = 1.072491980815252
$ws.Range("M23").Value = 1.081063086947769
$ws.Range("N23").Value = 1.072943754514889
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.067664470652325
$ws.Range("D24").Value = 1.070567155608349
$ws.Range("E24").Value = 1.07143859736981
$ws.Range("F24").Value = 1.080156896163665
$ws.Range("I24").Value = 1.050214903868837
$ws.Range("J24").Value = 1.073648393459488
$ws.Range("K24").Value = 1.07382086199378
$ws.Range("L24").Value = 1.074689476784396
$ws.Range("M24").Value = 1.083379768076813
$ws.Range("N24").Value = 1.075173097327392
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070822430975431
$ws.Range("D25").Value = 1.073582068912281
$ws.Range("E25").Value = 1.074283280752341
$ws.Range("F25").Value = 1.08313853499963
$ws.Range("I25").Value = 1.051125365932945
$ws.Range("J25").Value = 1.076212292211472
$ws.Range("K25").Value = 1.076520409343925
$ws.Range("L25").Value = 1.077219584499847
$ws.Range("M25").Value = 1.086049364966487
$ws.Range("N25").Value = 1.077740637109688
